$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing odds values for rows 2-7 (per upstream odds refresh) ---
# Row 2
$ws.Range("F2").Value = 1.73
$ws.Range("H2").Value = 5.4
$ws.Range("I2").Value = 5.6
$ws.Range("L2").Value = 1.37
$ws.Range("N2").Value = 4.5
$ws.Range("U2").Value = 2.16
$ws.Range("AC2").Value = 9
$ws.Range("AG2").Value = 9.800000000000001
$ws.Range("AJ2").Value = 17.5
$ws.Range("AM2").Value = 980
$ws.Range("AN2").Value = 9.6

# Row 3
$ws.Range("J3").Value = 1.09
$ws.Range("K3").Value = 980
$ws.Range("M3").Value = 1.02
$ws.Range("N3").Value = 1.1
$ws.Range("O3").Value = 1.24
$ws.Range("P3").Value = 1.45
$ws.Range("Q3").Value = 1.25
$ws.Range("S3").Value = 1.05
$ws.Range("T3").Value = 1.04
$ws.Range("U3").Value = 1.04

# Row 4
$ws.Range("F4").Value = 1.35
$ws.Range("G4").Value = 1.42
$ws.Range("H4").Value = 13
$ws.Range("I4").Value = 17
$ws.Range("J4").Value = 4.4
$ws.Range("K4").Value = 5.1
$ws.Range("L4").Value = 1.47
$ws.Range("N4").Value = 2.94
$ws.Range("W4").Value = 3.35
$ws.Range("X4").Value = 12
$ws.Range("Y4").Value = 30
$ws.Range("AC4").Value = 12
$ws.Range("AD4").Value = 70
$ws.Range("AG4").Value = 12
$ws.Range("AJ4").Value = 11.5
$ws.Range("AK4").Value = 21
$ws.Range("AL4").Value = 190

# Row 5
$ws.Range("O5").Value = 1.29
$ws.Range("Q5").Value = 1.29

# Row 6
$ws.Range("F6").Value = 8.199999999999999
$ws.Range("G6").Value = 8.6
$ws.Range("L6").Value = 1.36
$ws.Range("N6").Value = 4.1
$ws.Range("O6").Value = 1.3
$ws.Range("P6").Value = 2.06
$ws.Range("Q6").Value = 1.91
$ws.Range("R6").Value = 1.41
$ws.Range("S6").Value = 3.25
$ws.Range("T6").Value = 2.1
$ws.Range("U6").Value = 1.84
$ws.Range("W6").Value = 1.13
$ws.Range("X6").Value = 16.5
$ws.Range("Y6").Value = 7.8
$ws.Range("Z6").Value = 8.199999999999999
$ws.Range("AA6").Value = 12.5
$ws.Range("AB6").Value = 25
$ws.Range("AC6").Value = 10.5
$ws.Range("AD6").Value = 9.800000000000001
$ws.Range("AF6").Value = 75
$ws.Range("AG6").Value = 30
$ws.Range("AJ6").Value = 300
$ws.Range("AK6").Value = 140
$ws.Range("AL6").Value = 130
$ws.Range("AM6").Value = 170
$ws.Range("AN6").Value = 180
$ws.Range("AO6").Value = 7.8

# Row 7
$ws.Range("Q7").Value = 1.34
$ws.Range("S7").Value = 1.34

# --- Append new match row 8: Honduras Liga Nacional, Real Espana vs CD Motagua ---
$ws.Range("A8").Value = "Honduras Liga Nacional"
$ws.Range("D8").Value = "Real Espana"
$ws.Range("E8").Value = "CD Motagua"

# Date/time-look-alike text must stay literal strings, not auto-convert to Excel serials,
# so enter with a leading quote (forces text) then strip the resulting cell format.
$ws.Range("B8").Value = "'2025-12-23"
$ws.Range("B8").ClearFormats()
$ws.Range("C8").Value = "'22:00:00"
$ws.Range("C8").ClearFormats()

$ws.Range("F8").Value = 1.04
$ws.Range("G8").Value = 1000
$ws.Range("H8").Value = 1.04
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 1.03
$ws.Range("K8").Value = 950
$ws.Range("L8").Value = 1.01
$ws.Range("M8").Value = 1.01
$ws.Range("N8").Value = 1.34
$ws.Range("O8").Value = 1.01
$ws.Range("P8").Value = 1.34
$ws.Range("Q8").Value = 1.36
$ws.Range("R8").Value = 1.18
$ws.Range("S8").Value = 1.37
$ws.Range("T8").Value = 1.04
$ws.Range("U8").Value = 1.04
$ws.Range("V8").Value = 1.01
$ws.Range("W8").Value = 1.01
$ws.Range("X8").Value = 1000
$ws.Range("Y8").Value = 1000
$ws.Range("Z8").Value = 1000
$ws.Range("AA8").Value = 1000
$ws.Range("AB8").Value = 1000
$ws.Range("AC8").Value = 1000
$ws.Range("AD8").Value = 1000
$ws.Range("AE8").Value = 1000
$ws.Range("AF8").Value = 1000
$ws.Range("AG8").Value = 1000
$ws.Range("AH8").Value = 1000
$ws.Range("AI8").Value = 1000
$ws.Range("AJ8").Value = 1000
$ws.Range("AK8").Value = 1000
$ws.Range("AL8").Value = 1000
$ws.Range("AM8").Value = 1000
$ws.Range("AN8").Value = 1000
$ws.Range("AO8").Value = 1000
